$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold + border + centered) from G1 onto the new H1
# header cell so the new column matches the existing header formatting
# exactly (reuses the same style index rather than creating a new one).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
